# Generate Report for Handoff
# Updates the localization-status report: status moves from "In Translation"
# to "Ready for handoff", and the related timestamps are refreshed to the
# handoff generation time. Column widths are widened to fit the new,
# longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: per-language status + the overall "Latest HO Xliff
# Generate Date" timestamp.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-07 09:01:37"

# zh-cn detail sheet: status + its own "Latest Handoff Datetime".
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-07 09:01:29"

# de-de detail sheet: status + its own "Latest Handoff Datetime"
# (matches the Overview's generate date for this run).
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-07 09:01:37"

# Widen the status columns so the longer "Ready for handoff" text fits.
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
